# Auto-generated: update F-column ('想去人数' / want-to-go count) values
# across all four sheets to match refreshed data snapshot.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 112
$ws.Cells.Item(3, 6).Value = 1272
$ws.Cells.Item(4, 6).Value = 904
$ws.Cells.Item(5, 6).Value = 945
$ws.Cells.Item(6, 6).Value = 1695
$ws.Cells.Item(7, 6).Value = 369
$ws.Cells.Item(8, 6).Value = 1135
$ws.Cells.Item(9, 6).Value = 46
$ws.Cells.Item(10, 6).Value = 2
$ws.Cells.Item(11, 6).Value = 101
$ws.Cells.Item(12, 6).Value = 255
$ws.Cells.Item(13, 6).Value = 20
$ws.Cells.Item(14, 6).Value = 76
$ws.Cells.Item(15, 6).Value = 623
$ws.Cells.Item(16, 6).Value = 125
$ws.Cells.Item(17, 6).Value = 77
$ws.Cells.Item(19, 6).Value = 115
$ws.Cells.Item(20, 6).Value = 314
$ws.Cells.Item(21, 6).Value = 87
$ws.Cells.Item(22, 6).Value = 638
$ws.Cells.Item(23, 6).Value = 7
$ws.Cells.Item(24, 6).Value = 621
$ws.Cells.Item(25, 6).Value = 121
$ws.Cells.Item(26, 6).Value = 27
$ws.Cells.Item(28, 6).Value = 290
$ws.Cells.Item(29, 6).Value = 78
$ws.Cells.Item(30, 6).Value = 20
$ws.Cells.Item(31, 6).Value = 239
$ws.Cells.Item(32, 6).Value = 5
$ws.Cells.Item(34, 6).Value = 392

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(4, 6).Value = 305
$ws.Cells.Item(7, 6).Value = 230
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(10, 6).Value = 611
$ws.Cells.Item(11, 6).Value = 105

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 293

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 293
$ws.Cells.Item(3, 6).Value = 112
$ws.Cells.Item(4, 6).Value = 1272
$ws.Cells.Item(5, 6).Value = 904
$ws.Cells.Item(6, 6).Value = 945
$ws.Cells.Item(7, 6).Value = 1695
$ws.Cells.Item(8, 6).Value = 369
$ws.Cells.Item(9, 6).Value = 1135
$ws.Cells.Item(10, 6).Value = 46
$ws.Cells.Item(12, 6).Value = 2
$ws.Cells.Item(13, 6).Value = 101
$ws.Cells.Item(14, 6).Value = 255
$ws.Cells.Item(15, 6).Value = 20
$ws.Cells.Item(16, 6).Value = 76
$ws.Cells.Item(17, 6).Value = 623
$ws.Cells.Item(18, 6).Value = 125
$ws.Cells.Item(19, 6).Value = 77
$ws.Cells.Item(22, 6).Value = 305
$ws.Cells.Item(23, 6).Value = 115
$ws.Cells.Item(25, 6).Value = 314
$ws.Cells.Item(27, 6).Value = 230
$ws.Cells.Item(28, 6).Value = 230
$ws.Cells.Item(29, 6).Value = 87
$ws.Cells.Item(30, 6).Value = 638
$ws.Cells.Item(31, 6).Value = 7
$ws.Cells.Item(32, 6).Value = 621
$ws.Cells.Item(33, 6).Value = 121
$ws.Cells.Item(34, 6).Value = 27
$ws.Cells.Item(36, 6).Value = 290
$ws.Cells.Item(38, 6).Value = 1
$ws.Cells.Item(39, 6).Value = 78
$ws.Cells.Item(40, 6).Value = 20
$ws.Cells.Item(41, 6).Value = 239
$ws.Cells.Item(42, 6).Value = 611
$ws.Cells.Item(43, 6).Value = 105
$ws.Cells.Item(44, 6).Value = 105
$ws.Cells.Item(45, 6).Value = 5
$ws.Cells.Item(48, 6).Value = 392
